$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellRuns {
    param($cell, [string[]]$texts)

    $runsXml = ""
    foreach ($txt in $texts) {
        $runsXml += "<w:r><w:t>$txt</w:t></w:r>"
    }

    $xml = '<?xml version="1.0" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$cell.Range.InsertXML($xml)
}

# Row 2 -> "10X10"
Set-CellRuns $t.Cell(2, 2) @("2.843")
Set-CellRuns $t.Cell(2, 3) @("11.487")
Set-CellRuns $t.Cell(2, 4) @("42.125")
Set-CellRuns $t.Cell(2, 5) @("168", ".875")

# Row 3 -> "25X25"
Set-CellRuns $t.Cell(3, 2) @("41.812")
Set-CellRuns $t.Cell(3, 3) @("168.486")
Set-CellRuns $t.Cell(3, 4) @("66", "7.125")
Set-CellRuns $t.Cell(3, 5) @("2675.500")

# Row 4 -> "40X40"
Set-CellRuns $t.Cell(4, 2) @("168.750")
Set-CellRuns $t.Cell(4, 3) @("675.125")
Set-CellRuns $t.Cell(4, 4) @("2700.437")
Set-CellRuns $t.Cell(4, 5) @("10802.843")
